$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Progress Tracker")

# --- Row 33: NOT STARTED/Itertools-Infinie-Iterators -> COMPLETE/Itertools-Infinite-Iterators, dated 12/27/2020
$ws.Range("A33").Value = "COMPLETE"
$ws.Range("E33").Value = "Itertools - Infinite Iterators"
$ws.Range("F33").Value = "Correct Link: https://github.com/nyu-cds/python-itertools/tree/master/_episodes"

# --- Row 34: NOT STARTED/Itertools-Terminating-Iterators -> COMPLETE, dated 12/29/2020
$ws.Range("A34").Value = "COMPLETE"
$ws.Range("E34").Value = "Itertools - Terminating Iterators"
$ws.Range("F34").Value = "Correct Link: https://github.com/nyu-cds/python-itertools/tree/master/_episodes"

# --- Row 35: NOT STARTED/Itertools-Combinatoric-Generators -> COMPLETE, dated 12/29/2020
$ws.Range("A35").Value = "COMPLETE"
$ws.Range("E35").Value = "Itertools - Combinatoric Generators"
$ws.Range("F35").Value = "Correct Link: https://github.com/nyu-cds/python-itertools/tree/master/_episodes"

# --- Row 36: NOT STARTED/Assignment 5 -> SKIPPED/Assignment 5 - Part 1, dated 12/29/2020
$ws.Range("A36").Value = "SKIPPED"
$ws.Range("E36").Value = "Assignment 5 - Part 1"
$ws.Range("F36").Value = "Already used itertools in assignment 3"

# Apply the existing date-formatted style (as used by B2, etc.) to B33:B36, then set their date values
$ws.Range("B2").Copy()
$ws.Range("B33:B36").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B33").Value = 44192
$ws.Range("B34").Value = 44194
$ws.Range("B35").Value = 44194
$ws.Range("B36").Value = 44194

# --- Row 37: was entirely blank -> new row: IN PROGRESS / Exercise / Assignemnt 5 - Part 2
$ws.Range("A37").Value = "IN PROGRESS"
$ws.Range("C37").Value = 5
$ws.Range("D37").Value = "Exercise"
$ws.Range("E37").Value = "Assignemnt 5 - Part 2"

# --- Update the active selection in the sheet view to B37 (matches the author's final cursor position)
$ws.Range("B37").Select()
